# Apply cryptos list update (prices/volumes refreshed; two coin rows swapped)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.108.62"
$ws.Range("E2").Value = "  -1.52%  "

$ws.Range("D3").Value = "1.554.08"
$ws.Range("E3").Value = "  -0.89%  "

$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  -0.16%  "

$ws.Range("D5").Value = "'1.000"
$ws.Range("E5").Value = "  -0.14%  "

$ws.Range("D6").Value = "'287.14"

$ws.Range("D7").Value = "'0.3833"
$ws.Range("E7").Value = "  +2.88%  "

$ws.Range("D8").Value = "'0.3280"
$ws.Range("E8").Value = "  -1.21%  "

$ws.Range("D9").Value = "'43.68"
$ws.Range("E9").Value = "  -9.42%  "

$ws.Range("D10").Value = "'1.129"
$ws.Range("E10").Value = "  -0.40%  "

$ws.Range("D11").Value = "'0.07366"
$ws.Range("E11").Value = "  -1.44%  "

$ws.Range("D12").Value = "'1.000"
$ws.Range("E12").Value = "  -0.17%  "

$ws.Range("D13").Value = "'20.10"
$ws.Range("E13").Value = "  -2.74%  "

$ws.Range("D14").Value = "'5.797"
$ws.Range("E14").Value = "  -2.29%  "

$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "1.590.90"
$ws.Range("E15").Value = "  +1.75%  "

$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").Value = "'6.755"
$ws.Range("E16").Value = "  -2.28%  "

$ws.Range("D17").Value = "'0.00001089"
$ws.Range("E17").Value = "  -2.39%  "

$ws.Range("E18").Value = "  -1.72%  "

$ws.Range("D19").Value = "'85.83"
$ws.Range("E19").Value = "  -2.22%  "

$ws.Range("E20").Value = "  -0.16%  "

$ws.Range("D21").Value = "'6.370"
$ws.Range("E21").Value = "  +0.36%  "

$ws.Range("D22").Value = "'16.04"
$ws.Range("E22").Value = "  -2.44%  "

$ws.Range("D23").Value = "'11.68"
$ws.Range("E23").Value = "  -3.17%  "

$ws.Range("D24").Value = "22.113.35"
$ws.Range("E24").Value = "  -1.48%  "

$ws.Range("D25").Value = "'2.300"
$ws.Range("E25").Value = "  -3.55%  "

$ws.Range("E26").Value = "  -2.48%  "

$ws.Range("D27").Value = "'150.65"
$ws.Range("E27").Value = "  -1.39%  "

$ws.Range("D28").Value = "'19.16"
$ws.Range("E28").Value = "  -2.63%  "

$ws.Range("D29").Value = "'4.933"
$ws.Range("E29").Value = "  -1.61%  "

$ws.Range("B30").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C30").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D30").Value = "1.760.49"
$ws.Range("E30").Value = "  +1.17%  "

$ws.Range("B31").Value = "BitcoinCash"
$ws.Range("C31").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D31").Value = "'121.57"
$ws.Range("E31").Value = "  -1.99%  "

$ws.Range("D32").Value = "'1.080"
$ws.Range("E32").Value = "  +2.34%  "

$ws.Range("D33").Value = "'5.887"
$ws.Range("E33").Value = "  -3.89%  "

$ws.Range("D34").Value = "'1.901"
$ws.Range("E34").Value = "  -5.41%  "

$ws.Range("D35").Value = "'0.08236"
$ws.Range("E35").Value = "  -0.76%  "

$ws.Range("D36").Value = "'9.246"
$ws.Range("E36").Value = "  -5.19%  "

$ws.Range("D37").Value = "'0.06299"
$ws.Range("E37").Value = "  -1.52%  "

$ws.Range("D38").Value = "'0.02317"
$ws.Range("E38").Value = "  -5.72%  "

$ws.Range("D39").Value = "'5.270"
$ws.Range("E39").Value = "  -1.94%  "

$ws.Range("D40").Value = "'0.2151"
$ws.Range("E40").Value = "  -5.34%  "

$ws.Range("D41").Value = "'1.231"
$ws.Range("E41").Value = "  -4.32%  "

$ws.Range("D42").Value = "'11.02"
$ws.Range("E42").Value = "  -2.30%  "

$ws.Range("E43").Value = "  -0.15%  "

$ws.Range("D44").Value = "'0.6011"
$ws.Range("E44").Value = "  -4.43%  "

$ws.Range("D45").Value = "'13.63"
$ws.Range("E45").Value = "  -1.50%  "

$ws.Range("E46").Value = "  -1.15%  "

$ws.Range("D47").Value = "'0.5823"
$ws.Range("E47").Value = "  -5.19%  "

$ws.Range("D48").Value = "'1.971"
$ws.Range("E48").Value = "  -3.77%  "

$ws.Range("D49").Value = "'121.87"
$ws.Range("E49").Value = "  -2.92%  "

$ws.Range("D50").Value = "'1.173"
$ws.Range("E50").Value = "  -3.03%  "

$ws.Range("D51").Value = "'0.07023"
$ws.Range("E51").Value = "  -2.75%  "
